$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell values for the updated crypto price/volume snapshot.
# Values are plain numeric-looking strings (prices) or percentage strings
# (volume change); every target cell in the source workbook is stored as
# text (inlineStr), so each write is prefixed with a literal leading
# apostrophe to force text entry (preventing Excel's automatic
# number/date inference) and then the cell style is reset to "Normal" so
# no stray quote-prefix formatting is left on the cell itself.
$updates = [ordered]@{
    "D2" = "30.517.67"
    "E2" = "  -0.39%  "
    "D3" = "1.874.29"
    "E3" = "  -0.95%  "
    "D4" = "1.000"
    "E4" = "  -0.06%  "
    "D5" = "236.23"
    "E5" = "  -3.33%  "
    "D6" = "0.9994"
    "E6" = "  -0.11%  "
    "D7" = "0.4866"
    "E7" = "  -1.92%  "
    "D8" = "0.2892"
    "E8" = "  -2.29%  "
    "D9" = "0.06666"
    "E9" = "  -2.09%  "
    "D10" = "1.871.92"
    "E10" = "  -1.07%  "
    "D11" = "16.58"
    "E11" = "  -3.06%  "
    "E12" = "  -1.44%  "
    "D13" = "89.43"
    "E13" = "  -1.91%  "
    "D14" = "4.998"
    "E14" = "  -2.03%  "
    "D15" = "0.6531"
    "E15" = "  -3.19%  "
    "D16" = "30.466.36"
    "E16" = "  -0.53%  "
    "D17" = "0.000007821"
    "E17" = "  -1.31%  "
    "D18" = "0.9996"
    "E18" = "  -0.06%  "
    "D19" = "13.00"
    "E19" = "  -2.06%  "
    "D20" = "2.113.06"
    "E20" = "  -1.03%  "
    "D21" = "1.000"
    "E21" = "  -0.07%  "
    "D22" = "213.06"
    "E22" = "  +18.96%  "
    "D23" = "4.729"
    "E23" = "  -2.74%  "
    "D24" = "6.125"
    "E24" = "  +1.12%  "
    "D25" = "9.366"
    "E25" = "  +0.65%  "
    "D26" = "156.17"
    "E26" = "  +1.32%  "
    "D27" = "19.07"
    "E27" = "  +1.56%  "
    "E28" = "  -5.17%  "
    "D29" = "1.412"
    "E29" = "  +1.72%  "
    "D30" = "4.255"
    "E30" = "  -1.72%  "
    "E31" = "  +1.26%  "
    "D32" = "3.916"
    "E32" = "  -3.03%  "
    "D33" = "0.05119"
    "E33" = "  -1.70%  "
    "D34" = "0.7246"
    "E34" = "  -1.63%  "
    "D35" = "1.077"
    "E35" = "  -4.92%  "
    "D36" = "2.687"
    "E36" = "  +0.29%  "
    "D37" = "0.01814"
    "E37" = "  -3.24%  "
    "D38" = "2.656"
    "E38" = "  -1.65%  "
    "D39" = "0.9191"
    "E39" = "  -1.75%  "
    "D40" = "2.044"
    "E40" = "  -5.76%  "
    "D41" = "0.4404"
    "E41" = "  +0.97%  "
    "D42" = "104.49"
    "E42" = "  -1.43%  "
    "D43" = "5.737"
    "E43" = "  -1.27%  "
    "E44" = "  -0.66%  "
    "E45" = "  -1.87%  "
    "D46" = "7.324"
    "E46" = "  -4.23%  "
    "B47" = "Cronos"
    "C47" = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
    "D47" = "0.05828"
    "E47" = "  -0.29%  "
    "B48" = "Decentraland"
    "C48" = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
    "D48" = "0.4008"
    "E48" = "  +2.90%  "
    "D49" = "8.605"
    "E49" = "  +0.79%  "
    "E50" = "  +2.08%  "
    "D51" = "33.21"
    "E51" = "  -0.51%  "
}

foreach ($addr in $updates.Keys) {
    $value = $updates[$addr]
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}
